# Apply updated cryptocurrency price/volume data to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

Set-TextValue "D2" '42.701.91'
Set-TextValue "E2" '  -7.52%  '

Set-TextValue "D3" '2.541.28'
Set-TextValue "E3" '  -4.14%  '

Set-TextValue "E4" '  +0.06%  '

Set-TextValue "D5" '299.13'
Set-TextValue "E5" '  -3.74%  '

Set-TextValue "D6" '92.21'
Set-TextValue "E6" '  -6.75%  '

Set-TextValue "E7" '  -3.98%  '

Set-TextValue "E8" '  +0.01%  '

Set-TextValue "D9" '0.550'
Set-TextValue "E9" '  -5.69%  '

Set-TextValue "D10" '35.80'
Set-TextValue "E10" '  -7.99%  '

Set-TextValue "D11" '0.0800'
Set-TextValue "E11" '  -5.59%  '

Set-TextValue "D12" '7.62'
Set-TextValue "E12" '  -5.96%  '

Set-TextValue "D13" '0.113'
Set-TextValue "E13" '  +4.89%  '

Set-TextValue "D14" '2.929.94'
Set-TextValue "E14" '  -4.04%  '

Set-TextValue "D15" '2.527.75'
Set-TextValue "E15" '  -4.43%  '

Set-TextValue "D16" '0.872'
Set-TextValue "E16" '  -5.67%  '

Set-TextValue "D17" '14.20'
Set-TextValue "E17" '  -5.06%  '

Set-TextValue "D18" '42.749.92'
Set-TextValue "E18" '  -7.49%  '

Set-TextValue "D19" '12.88'
Set-TextValue "E19" '  +0.61%  '

Set-TextValue "E20" '  -3.95%  '

Set-TextValue "D21" '6.55'
Set-TextValue "E21" '  -3.67%  '

Set-TextValue "D22" '71.26'
Set-TextValue "E22" '  -4.63%  '

Set-TextValue "D23" '256.40'
Set-TextValue "E23" '  -9.28%  '

Set-TextValue "D24" '2.91'
Set-TextValue "E24" '  -4.53%  '

Set-TextValue "D25" '29.21'
Set-TextValue "E25" '  -4.73%  '

Set-TextValue "D26" '2.12'
Set-TextValue "E26" '  -6.25%  '

Set-TextValue "E27" '  -0.38%  '

Set-TextValue "D28" '10.03'
Set-TextValue "E28" '  -5.19%  '

Set-TextValue "D29" '37.00'
Set-TextValue "E29" '  -4.51%  '

Set-TextValue "E30" '  -5.67%  '

Set-TextValue "E31" '  -5.35%  '

Set-TextValue "D32" '152.40'
Set-TextValue "E32" '  -3.08%  '

Set-TextValue "E33" '  -7.83%  '

Set-TextValue "E34" '  -2.55%  '

Set-TextValue "E35" '  -9.90%  '

Set-TextValue "D36" '0.0793'
Set-TextValue "E36" '  -6.05%  '

Set-TextValue "E37" '  -7.73%  '

Set-TextValue "D38" '0.119'
Set-TextValue "E38" '  -3.54%  '

Set-TextValue "B39" 'Celestia'
Set-TextValue "C39" 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D39" '17.03'
Set-TextValue "E39" '  +7.49%  '

Set-TextValue "B40" 'EnergySwap'
Set-TextValue "C40" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D40" '23.95'
Set-TextValue "E40" '  +4.26%  '

Set-TextValue "E41" '  -5.93%  '

Set-TextValue "D42" '3.87'
Set-TextValue "E42" '  -4.29%  '

Set-TextValue "E43" '  -5.46%  '

Set-TextValue "D44" '2.079.37'
Set-TextValue "E44" '  -3.02%  '

Set-TextValue "D45" '0.999'
Set-TextValue "E45" '  -0.05%  '

Set-TextValue "D46" '9.13'
Set-TextValue "E46" '  -0.43%  '

Set-TextValue "B47" 'ApeXProtocol'
Set-TextValue "C47" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D47" '1.62'
Set-TextValue "E47" '  +3.10%  '

Set-TextValue "B48" 'BitcoinSV'
Set-TextValue "C48" 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue "D48" '84.39'
Set-TextValue "E48" '  -10.36%  '

Set-TextValue "D49" '2.787.15'
Set-TextValue "E49" '  -3.82%  '

Set-TextValue "D50" '103.96'
Set-TextValue "E50" '  -6.10%  '

Set-TextValue "D51" '1.66'
Set-TextValue "E51" '  -5.25%  '

Write-Host "Updated cryptos list"
